# The RANK() calls in columns F, K, P of the metro_budget sheet were
# ranking descending (the default) when the author intended ascending
# order. Add the third RANK argument (1 = ascending) to every RANK call
# in those columns. Because F3:F52 / K3:K52 / P3:P52 are shared-formula
# groups, updating the anchor formula (row 3) updates every cell in the
# group; row 2 holds its own (non-shared) formula and is updated
# separately.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# Row 2 (non-shared formulas)
$ws.Range("F2").Formula = '=IFERROR(RANK(E2,E$2:E$52,1),"N/A")'
$ws.Range("K2").Formula = '=IFERROR(RANK(J2,J$2:J$52,1),"N/A")'
$ws.Range("P2").Formula = '=IFERROR(RANK(O2,O$2:O$52,1),"N/A")'

# Rows 3:52 (shared formula groups anchored at row 3)
$ws.Range("F3:F52").Formula = '=IFERROR(RANK(E3,E$2:E$52,1),"N/A")'
$ws.Range("K3:K52").Formula = '=IFERROR(RANK(J3,J$2:J$52,1),"N/A")'
$ws.Range("P3:P52").Formula = '=IFERROR(RANK(O3,O$2:O$52,1),"N/A")'

# Restore the view state recorded in the saved workbook: the active
# sheet's scroll position and the selected cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$ws.Range("M11").Select()
